$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 - Create Tag Functionality
$ws.Range("A19").Value = 43504
$ws.Range("A19").NumberFormat = $ws.Range("A18").NumberFormat
$ws.Range("B19").Value = 0.25
$ws.Range("C19").Value = "Implementation"
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = "Create Tag Functionality"
$ws.Range("F19").Value = "Once you" + [char]0x2019 + "ve done one" + [char]0x2026

# Row 20 - Create Merchant Functionality
$ws.Range("A20").Value = 43504
$ws.Range("A20").NumberFormat = $ws.Range("A18").NumberFormat
$ws.Range("B20").Value = 0.25
$ws.Range("C20").Value = "Implementation"
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = "Create Merchant Functionality"
$ws.Range("F20").Value = "Once you" + [char]0x2019 + "ve done one" + [char]0x2026

# Match existing style of last data row (A18:F18) for the new rows
$ws.Range("A18:F18").Copy()
$ws.Range("A19:F20").PasteSpecial(-4122)

# Re-set the selection as in the target sheet
$ws.Range("D12").Select()
